$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "69 x 94" + [char]11 + "  9    4" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "9|    |"
$t.Cell(1,2).Range.Text = "26 x 93" + [char]11 + "  9    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "6|    |"
$t.Cell(1,3).Range.Text = "94 x 85" + [char]11 + "  8    5" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "4|    |"
$t.Cell(2,1).Range.Text = "56 x 28" + [char]11 + "  2    8" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "6|    |"
$t.Cell(2,2).Range.Text = "67 x 19" + [char]11 + "  1    9" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "7|    |"
$t.Cell(2,3).Range.Text = "47 x 90" + [char]11 + "  9    0" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "7|    |"
$t.Cell(3,1).Range.Text = "92 x 53" + [char]11 + "  5    3" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "2|    |"
$t.Cell(3,2).Range.Text = "65 x 10" + [char]11 + "  1    0" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "5|    |"
$t.Cell(3,3).Range.Text = "97 x 72" + [char]11 + "  7    2" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
$t.Cell(4,1).Range.Text = "43 x 86" + [char]11 + "  8    6" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "3|    |"
$t.Cell(4,2).Range.Text = "17 x 94" + [char]11 + "  9    4" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "7|    |"
$t.Cell(4,3).Range.Text = "73 x 71" + [char]11 + "  7    1" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "3|    |"
$t.Cell(5,1).Range.Text = "78 x 36" + [char]11 + "  3    6" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "8|    |"
$t.Cell(5,2).Range.Text = "28 x 24" + [char]11 + "  2    4" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "8|    |"
$t.Cell(5,3).Range.Text = "23 x 22" + [char]11 + "  2    2" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "3|    |"
